$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: nueva venta (corte de caja)
# A4 looks like a date ("2023-10-14"); force it to stay plain text like the
# other date cells in the sheet (A2/A3) instead of being auto-converted to a
# date serial number by Excel's smart-entry parsing.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2023-10-14"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "cono sencillo, cono sencillo"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = "tarjeta"
